# Fix priority room by CS, campus, backup
#
# The "room" sheet lists classrooms with their Building assignment in
# column B. Rows 16-33 (RoomID ClassRoom1_200 .. ClassRoom6_300) were all
# tagged "campus" - repurpose them as the "backup" priority pool instead.
#
# Also update the sheet's current selection to reflect where the author
# left off reviewing the change (G21), and drop the stale scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("room")

# Re-tag the 18 "campus" rooms (rows 16-33) as "backup".
$ws.Range("B16:B33").Value = "backup"

# Move/replace the selection to G21 (also clears the old topLeftCell scroll
# anchor and the previous F39 selection).
$ws.Range("G21").Select()
